# Update the "Förändrad" (Changed) date column (C) from serial 45177 to 45178
# for all data rows (2 through 189) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 189
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)   # Column C
    if ($cell.Value2 -eq 45177) {
        $cell.Value = 45178
    }
}
